$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (Price / Volume) to Text format so values such as
# "1.001" or "30.478.47" are not auto-converted to numbers/dates by Excel,
# matching the inline-string cells produced by the scraper.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '30.478.47'
$ws.Range("E2").Value = '  -0.47%  '

# Row 3
$ws.Range("D3").Value = '1.889.01'
$ws.Range("E3").Value = '  +0.73%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").Value = '243.44'
$ws.Range("E5").Value = '  -1.71%  '

# Row 6
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.07%  '

# Row 7
$ws.Range("D7").Value = '0.4706'
$ws.Range("E7").Value = '  -0.52%  '

# Row 8
$ws.Range("D8").Value = '0.2897'
$ws.Range("E8").Value = '  -0.39%  '

# Row 9
$ws.Range("D9").Value = '0.06492'
$ws.Range("E9").Value = '  +0.16%  '

# Row 10
$ws.Range("D10").Value = '22.19'
$ws.Range("E10").Value = '  +0.27%  '

# Row 11
$ws.Range("D11").Value = '0.07750'
$ws.Range("E11").Value = '  +0.55%  '

# Row 12
$ws.Range("D12").Value = '1.888.92'
$ws.Range("E12").Value = '  +0.74%  '

# Row 13
$ws.Range("D13").Value = '95.65'
$ws.Range("E13").Value = '  -0.71%  '

# Row 14
$ws.Range("D14").Value = '0.7252'
$ws.Range("E14").Value = '  -1.88%  '

# Row 15
$ws.Range("D15").Value = '5.188'
$ws.Range("E15").Value = '  +0.48%  '

# Row 16
$ws.Range("D16").Value = '281.24'
$ws.Range("E16").Value = '  +2.83%  '

# Row 17
$ws.Range("D17").Value = '30.477.36'
$ws.Range("E17").Value = '  -0.58%  '

# Row 18
$ws.Range("D18").Value = '13.05'
$ws.Range("E18").Value = '  -2.10%  '

# Row 19
$ws.Range("E19").Value = '  +0.05%  '

# Row 20
$ws.Range("D20").Value = '0.000007465'
$ws.Range("E20").Value = '  -0.65%  '

# Row 21
$ws.Range("D21").Value = '2.137.11'
$ws.Range("E21").Value = '  +0.99%  '

# Row 22
$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.01%  '

# Row 23
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '5.281'
$ws.Range("E23").Value = '  +0.24%  '

# Row 24
$ws.Range("D24").Value = '6.261'
$ws.Range("E24").Value = '  +1.10%  '

# Row 25
$ws.Range("D25").Value = '163.89'
$ws.Range("E25").Value = '  -0.05%  '

# Row 26
$ws.Range("D26").Value = '9.069'
$ws.Range("E26").Value = '  -1.57%  '

# Row 27
$ws.Range("D27").Value = '18.89'
$ws.Range("E27").Value = '  +0.80%  '

# Row 28
$ws.Range("D28").Value = '1.895'
$ws.Range("E28").Value = '  -0.90%  '

# Row 29
$ws.Range("D29").Value = '0.09734'
$ws.Range("E29").Value = '  -2.88%  '

# Row 30
$ws.Range("D30").Value = '1.330'
$ws.Range("E30").Value = '  -1.16%  '

# Row 31
$ws.Range("E31").Value = '  -2.71%  '

# Row 32
$ws.Range("D32").Value = '4.277'
$ws.Range("E32").Value = '  +0.02%  '

# Row 33
$ws.Range("D33").Value = '4.143'
$ws.Range("E33").Value = '  +1.10%  '

# Row 34
$ws.Range("D34").Value = '0.04852'
$ws.Range("E34").Value = '  +1.17%  '

# Row 35
$ws.Range("D35").Value = '1.126'
$ws.Range("E35").Value = '  +0.42%  '

# Row 36
$ws.Range("D36").Value = '0.6948'
$ws.Range("E36").Value = '  -0.10%  '

# Row 37
$ws.Range("D37").Value = '2.719'
$ws.Range("E37").Value = '  +0.04%  '

# Row 38
$ws.Range("D38").Value = '0.01890'
$ws.Range("E38").Value = '  +2.13%  '

# Row 39
$ws.Range("D39").Value = '2.820'
$ws.Range("E39").Value = '  +2.53%  '

# Row 40
$ws.Range("D40").Value = '75.31'
$ws.Range("E40").Value = '  +2.89%  '

# Row 41
$ws.Range("D41").Value = '6.211'

# Row 42
$ws.Range("D42").Value = '1.985'
$ws.Range("E42").Value = '  +0.69%  '

# Row 43
$ws.Range("D43").Value = '0.4243'
$ws.Range("E43").Value = '  +1.52%  '

# Row 44
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.00%  '

# Row 45
$ws.Range("D45").Value = '0.8247'
$ws.Range("E45").Value = '  -1.26%  '

# Row 46
$ws.Range("D46").Value = '101.39'
$ws.Range("E46").Value = '  -0.88%  '

# Row 47
$ws.Range("D47").Value = '9.567'
$ws.Range("E47").Value = '  +2.02%  '

# Row 48
$ws.Range("D48").Value = '6.955'
$ws.Range("E48").Value = '  -0.54%  '

# Row 49
$ws.Range("D49").Value = '35.08'
$ws.Range("E49").Value = '  -1.02%  '

# Row 50
$ws.Range("D50").Value = '911.37'
$ws.Range("E50").Value = '  -0.82%  '

# Row 51
$ws.Range("E51").Value = '  +1.77%  '

# Remove the temporary Text number formatting so the cells end up with the
# same (default/no) style as before, while keeping the values as text.
$ws.Range("D2:E51").ClearFormats()
